# The workbook's single data sheet is being renamed as part of unifying the
# DataNode / DataTable / Entity naming convention, and the active selection
# is moved to C38 (reflecting where editing left off before saving).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Property1")
$ws.Name = "DataNode"

$ws.Range("C38").Select()
